# Groepadressen ETS.xlsx - research files for device values
# - Remove the now-unused "Datatypes" lookup sheet (its Type/ID table was
#   folded directly into the "Data_type(ID)" note columns on "Objects").
# - Refresh the KNX datapoint-type codes shown in the "Data_type(ID)" columns
#   of the "Objects" sheet to the corrected values.
# - Leave the view scrolled back to the top, with the selection parked on A30.

$excel.DisplayAlerts = $false

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Objects")

# Drop the separate "Datatypes" sheet; its content is no longer needed as a
# standalone tab.
$wb.Worksheets.Item("Datatypes").Delete()

# --- Update Data_type(ID) codes on "Objects" ---

# "Verlichting" (Schakelen/Dimmen) block - rows 3-11, column G
$ws.Range("G3:G11").Value = "1,1,15,2,2"

# "Blinds" block - rows 26-31, column G
$ws.Range("G26:G31").Value = "1,1,2,2"

# "Ventilation" (Bedroom/Livingroom/Bathroom) block - rows 19-21, column F
$ws.Range("F19:F21").Value = "4,2,4,-"

# "MTP dummy" (Color scene / Temperature color) block - rows 40-42, column H
$ws.Range("H40:H42").Value = "1,2,2,2,4"

# Reset the view: scroll back so column A is visible and park the selection.
$ws.Activate()
$ws.Range("A30").Select()
